$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.005", "24.505.32") that must
# stay as text rather than being auto-coerced into numbers by Excel input
# parsing. Force the Text number format first, matching how these cells are
# already stored (t="inlineStr") in the source workbook.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '24.505.32'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '1.656.64'
$ws.Range('E3').Value = '  -2.71%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = '307.06'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').Value = '0.3610'
$ws.Range('E7').Value = '  -3.43%  '
$ws.Range('D8').Value = '47.56'
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('D9').Value = '0.3238'
$ws.Range('E9').Value = '  -5.82%  '
$ws.Range('D10').Value = '1.116'
$ws.Range('E10').Value = '  -5.64%  '
$ws.Range('D11').Value = '0.06959'
$ws.Range('E11').Value = '  -6.76%  '
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '5.861'
$ws.Range('E13').Value = '  -5.99%  '
$ws.Range('D14').Value = '19.36'
$ws.Range('E14').Value = '  -7.33%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.648.19'
$ws.Range('E15').Value = '  -3.34%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '6.537'
$ws.Range('E16').Value = '  -5.61%  '
$ws.Range('D17').Value = '0.00001042'
$ws.Range('E17').Value = '  -7.09%  '
$ws.Range('D18').Value = '0.06525'
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('D19').Value = '0.9996'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '76.24'
$ws.Range('E20').Value = '  -8.83%  '
$ws.Range('D21').Value = '5.895'
$ws.Range('E21').Value = '  -6.82%  '
$ws.Range('D22').Value = '15.59'
$ws.Range('E22').Value = '  -8.78%  '
$ws.Range('D23').Value = '12.60'
$ws.Range('E23').Value = '  -4.59%  '
$ws.Range('D24').Value = '24.498.50'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = '2.464'
$ws.Range('E25').Value = '  +2.32%  '
$ws.Range('D26').Value = '2.286'
$ws.Range('E26').Value = '  -17.04%  '
$ws.Range('D27').Value = '146.79'
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('D28').Value = '18.38'
$ws.Range('E28').Value = '  -8.57%  '
$ws.Range('D29').Value = '1.837.14'
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('D30').Value = '1.188'
$ws.Range('E30').Value = '  +0.60%  '
$ws.Range('D31').Value = '123.63'
$ws.Range('E31').Value = '  -5.66%  '
$ws.Range('D32').Value = '4.080'
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('D33').Value = '5.632'
$ws.Range('E33').Value = '  -16.63%  '
$ws.Range('E34').Value = '  -4.43%  '
$ws.Range('D35').Value = '0.08348'
$ws.Range('E35').Value = '  -5.07%  '
$ws.Range('D36').Value = '12.31'
$ws.Range('E36').Value = '  -9.82%  '
$ws.Range('D37').Value = '5.182'
$ws.Range('E37').Value = '  -6.05%  '
$ws.Range('D38').Value = '0.06036'
$ws.Range('E38').Value = '  -7.20%  '
$ws.Range('D39').Value = '0.02186'
$ws.Range('E39').Value = '  -8.07%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.201'
$ws.Range('E40').Value = '  -5.61%  '
$ws.Range('D41').Value = '0.2047'
$ws.Range('E41').Value = '  -7.71%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '8.184'
$ws.Range('E42').Value = '  -8.25%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').Value = '0.5875'
$ws.Range('E44').Value = '  -8.18%  '
$ws.Range('D45').Value = '3.733'
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').Value = '12.68'
$ws.Range('E46').Value = '  -8.37%  '
$ws.Range('D47').Value = '0.5568'
$ws.Range('E47').Value = '  -8.49%  '
$ws.Range('D48').Value = '122.12'
$ws.Range('E48').Value = '  -5.52%  '
$ws.Range('D49').Value = '1.933'
$ws.Range('E49').Value = '  -8.47%  '
$ws.Range('D50').Value = '0.06915'
$ws.Range('E50').Value = '  -4.81%  '
$ws.Range('D51').Value = '73.95'
$ws.Range('E51').Value = '  -6.40%  '
